# This workbook's data rows (2-17) got shuffled: the values in columns
# D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen) and S (Precio $/Kg) move
# between rows while every other column stays put. Below is the mapping
# of destination row -> source row (the row whose original values end up
# in the destination row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 15
    3  = 8
    4  = 2
    5  = 3
    6  = 17
    7  = 14
    8  = 4
    9  = 13
    10 = 11
    11 = 16
    12 = 10
    13 = 6
    14 = 9
    15 = 7
    16 = 12
    17 = 5
}

$cols = @("D", "M", "N", "O", "P", "R", "S")

# First, snapshot the original values of every relevant cell so that
# overwriting a row doesn't clobber data still needed as a source for
# another row.
$snapshot = @{}
foreach ($row in 2..17) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Now write each destination row's cells using the snapshotted source
# row's values.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcValues[$col]
    }
}
